# edit device query setting
# Fill in start_date (column B) and expiry_date (column C) for the 12 sim_id
# data rows (rows 2-13) that were previously left blank, and move the
# active selection to the newly-populated expiry_date column (C2:C13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B -> start_date, Column C -> expiry_date
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = "26-3-2019"
    $ws.Cells.Item($r, 3).Value = "26-4-2019"
}

# Update the active selection to reflect the expiry_date column being reviewed.
$ws.Range("C2:C13").Select()
